# Updated cryptos list on Fri Jul 19 06:43:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# All target cells are formatted as Text first so numeric-looking
# strings (e.g. '0.427', '13.70', '43.00') keep their exact literal
# representation instead of being auto-coerced into numbers.
$targets = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","E9","E10","E11","D12","E12","D13","E13","E14","E15","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","E22","D23","E23","D24","E24","D25","E25","E26","D27","E27","E28","D29","E29","D30","E30","E31","E32","D33","E33","D34","E34","E35","E36","D37","E37","E38","B39","C39","D39","E39","B40","C40","D40","E40","B41","C41","D41","E41","D42","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","E49","D50","E50","E51")
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.295.55"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.439.57"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "573.88"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "164.87"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.441.74"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -5.41%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "0.427"
$ws.Range("E12").Value = "  -4.83%  "
$ws.Range("D13").Value = "4.034.93"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("E16").Value = "  -7.01%  "
$ws.Range("D17").Value = "64.352.66"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "3.509.75"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "13.70"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "379.74"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "71.60"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "6.12"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "23.03"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "7.15"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "0.863"
$ws.Range("E37").Value = "  +11.44%  "
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "26.28"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.822.83"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.0732"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").Value = "26.58"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "43.00"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "6.49"
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("D46").Value = "2.50"
$ws.Range("E46").Value = "  +9.94%  "
$ws.Range("D47").Value = "0.0309"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "335.63"
$ws.Range("E48").Value = "  +5.65%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "6.35"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  -3.04%  "
